$d = $word.ActiveDocument

# Update the worksheet date in the title paragraph.
$d.Content.Find.Execute("2023-08-04 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-05 Saturday", 2) | Out-Null

# Helper: replace text inside a single table cell, addressed by row/column.
# We re-derive a document-level Range (from the cell's own Start/End) before
# calling Find, since Find on a Range obtained straight from Cell.Range is not
# reliably bounded to that cell when the same source text also occurs elsewhere
# in the document (e.g. the two "73÷9=" cells below must get different results).
function Replace-CellText($table, $row, $col, $old, $new) {
    $cell = $table.Cell($row, $col)
    $cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
    $cellRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

# Update each division problem cell in the practice table.
$t = $d.Tables.Item(1)
Replace-CellText $t 1 1 "22÷7=" "31÷6="
Replace-CellText $t 1 2 "62÷2=" "72÷6="
Replace-CellText $t 1 3 "28÷9=" "88÷7="
Replace-CellText $t 1 4 "73÷9=" "30÷8="
Replace-CellText $t 1 5 "40÷3=" "25÷3="
Replace-CellText $t 5 1 "30÷4=" "69÷2="
Replace-CellText $t 5 2 "57÷2=" "78÷7="
Replace-CellText $t 5 3 "95÷5=" "94÷8="
Replace-CellText $t 5 4 "12÷6=" "56÷3="
Replace-CellText $t 5 5 "77÷4=" "78÷4="
Replace-CellText $t 9 1 "36÷2=" "88÷7="
Replace-CellText $t 9 2 "63÷7=" "13÷8="
Replace-CellText $t 9 3 "73÷9=" "95÷6="
Replace-CellText $t 9 4 "60÷5=" "98÷2="
Replace-CellText $t 9 5 "96÷4=" "59÷5="
Replace-CellText $t 13 1 "63÷9=" "85÷2="
Replace-CellText $t 13 2 "83÷2=" "96÷3="
Replace-CellText $t 13 3 "82÷9=" "87÷9="
Replace-CellText $t 13 4 "58÷4=" "79÷7="
Replace-CellText $t 13 5 "77÷3=" "96÷8="
Replace-CellText $t 17 1 "87÷2=" "20÷4="
Replace-CellText $t 17 2 "40÷4=" "39÷2="
Replace-CellText $t 17 3 "45÷9=" "31÷8="
Replace-CellText $t 17 4 "29÷7=" "94÷9="
Replace-CellText $t 17 5 "36÷6=" "90÷2="

Write-Output "done"
